$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "meaning"
$ws.Range("C1").Value = "synonym "
$ws.Range("D1").Value = "antonym"
$ws.Range("E1").Value = "usage"

$ws.Range("G5").Select()
